$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) with header + two data rows, matching the
# existing header formatting (same style as B1:G1 headers).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
